$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.424.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.296.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.93%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.44"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.286.10"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.614"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.06"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.842.12"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "17.99"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.304.35"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.63"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.394.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.960"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "422.48"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.60"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.05"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.39"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.87"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.78"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.56"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.50"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "583.96"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.30"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.80"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.42"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.97"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0740"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.360"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.067.17"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.79"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0401"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.45"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.06"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.21%  "
